$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 177 and row 178 (columns B:AC)
$row177 = New-Object 'object[,]' 1,28
$row177[0,0] = 6845311
$row177[0,1] = "Scotland Championship"
$row177[0,2] = "Scotland Championship"
$row177[0,3] = 45293.5
$row177[0,4] = "Dunfermline"
$row177[0,5] = "Raith"
$row177[0,6] = 1
$row177[0,7] = 2
$row177[0,8] = "A"
$row177[0,9] = 2.6
$row177[0,10] = 3.5
$row177[0,11] = 2.25
$row177[0,12] = 3.2
$row177[0,13] = 3.75
$row177[0,14] = 2
$row177[0,15] = 0.25
$row177[0,16] = 2
$row177[0,17] = 1.8
$row177[0,18] = 2.75
$row177[0,19] = 1.925
$row177[0,20] = 1.875
$row177[0,21] = -1
$row177[0,22] = -1
$row177[0,23] = 1
$row177[0,24] = -1
$row177[0,25] = 0.8
$row177[0,26] = 0.4625
$row177[0,27] = -0.5
$row178 = New-Object 'object[,]' 1,28
$row178[0,0] = 6845310
$row178[0,1] = "Scotland Championship"
$row178[0,2] = "Scotland Championship"
$row178[0,3] = 45293.5
$row178[0,4] = "Arbroath"
$row178[0,5] = "Dundee Utd"
$row178[0,6] = 0
$row178[0,7] = 3
$row178[0,8] = "A"
$row178[0,9] = 6
$row178[0,10] = 3.9
$row178[0,11] = 1.444
$row178[0,12] = 11
$row178[0,13] = 6
$row178[0,14] = 1.181
$row178[0,15] = 2
$row178[0,16] = 1.875
$row178[0,17] = 1.925
$row178[0,18] = 3.25
$row178[0,19] = 1.95
$row178[0,20] = 1.85
$row178[0,21] = -1
$row178[0,22] = -1
$row178[0,23] = 0.181
$row178[0,24] = -1
$row178[0,25] = 0.925
$row178[0,26] = -0.5
$row178[0,27] = 0.425
$ws.Range("B177:AC177").Value = $row178
$ws.Range("B178:AC178").Value = $row177

# Swap row 179 and row 181 (columns B:AC)
$row179 = New-Object 'object[,]' 1,28
$row179[0,0] = 6845314
$row179[0,1] = "Scotland Championship"
$row179[0,2] = "Scotland Championship"
$row179[0,3] = 45293.5
$row179[0,4] = "Airdrieonians"
$row179[0,5] = "Inverness CT"
$row179[0,6] = 2
$row179[0,7] = 0
$row179[0,8] = "H"
$row179[0,9] = 2.2
$row179[0,10] = 3.5
$row179[0,11] = 2.7
$row179[0,12] = 2.3
$row179[0,13] = 3.25
$row179[0,14] = 2.7
$row179[0,15] = 0
$row179[0,16] = 1.75
$row179[0,17] = 2.05
$row179[0,18] = 2.25
$row179[0,19] = 1.925
$row179[0,20] = 1.925
$row179[0,21] = 1.3
$row179[0,22] = -1
$row179[0,23] = -1
$row179[0,24] = 0.75
$row179[0,25] = -1
$row179[0,26] = -0.5
$row179[0,27] = 0.4625
$row181 = New-Object 'object[,]' 1,28
$row181[0,0] = 6845312
$row181[0,1] = "Scotland Championship"
$row181[0,2] = "Scotland Championship"
$row181[0,3] = 45293.5
$row181[0,4] = "Morton"
$row181[0,5] = "Ayr"
$row181[0,6] = 3
$row181[0,7] = 0
$row181[0,8] = "H"
$row181[0,9] = 2
$row181[0,10] = 3.4
$row181[0,11] = 3.1
$row181[0,12] = 2.1
$row181[0,13] = 3.4
$row181[0,14] = 2.9
$row181[0,15] = -0.25
$row181[0,16] = 1.9
$row181[0,17] = 1.9
$row181[0,18] = 2.5
$row181[0,19] = 1.925
$row181[0,20] = 1.875
$row181[0,21] = 1.1
$row181[0,22] = -1
$row181[0,23] = -1
$row181[0,24] = 0.8999999999999999
$row181[0,25] = -1
$row181[0,26] = 0.925
$row181[0,27] = -1
$ws.Range("B179:AC179").Value = $row181
$ws.Range("B181:AC181").Value = $row179

# Cyclic shift rows 211-215 (each row gets the following row's data; row 215 gets row 211's old data)
$row211_bg = New-Object 'object[,]' 1,6
$row211_bg[0,0] = 6853957
$row211_bg[0,1] = "Scotland Championship"
$row211_bg[0,2] = "Scotland Championship"
$row211_bg[0,3] = 45349.69791666666
$row211_bg[0,4] = "Raith"
$row211_bg[0,5] = "Morton"
$row211_kv = New-Object 'object[,]' 1,12
$row211_kv[0,0] = 1.8
$row211_kv[0,1] = 3.6
$row211_kv[0,2] = 3.8
$row211_kv[0,3] = 1.85
$row211_kv[0,4] = 3.6
$row211_kv[0,5] = 3.6
$row211_kv[0,6] = -0.5
$row211_kv[0,7] = 1.925
$row211_kv[0,8] = 1.925
$row211_kv[0,9] = 2.5
$row211_kv[0,10] = 1.875
$row211_kv[0,11] = 1.975
$row212_bg = New-Object 'object[,]' 1,6
$row212_bg[0,0] = 6845345
$row212_bg[0,1] = "Scotland Championship"
$row212_bg[0,2] = "Scotland Championship"
$row212_bg[0,3] = 45349.69791666666
$row212_bg[0,4] = "Ayr"
$row212_bg[0,5] = "Partick"
$row212_kv = New-Object 'object[,]' 1,12
$row212_kv[0,0] = 3
$row212_kv[0,1] = 3.5
$row212_kv[0,2] = 2.1
$row212_kv[0,3] = 2.875
$row212_kv[0,4] = 3.5
$row212_kv[0,5] = 2.15
$row212_kv[0,6] = 0.25
$row212_kv[0,7] = 1.875
$row212_kv[0,8] = 1.975
$row212_kv[0,9] = 2.75
$row212_kv[0,10] = 1.925
$row212_kv[0,11] = 1.925
$row213_bg = New-Object 'object[,]' 1,6
$row213_bg[0,0] = 6845346
$row213_bg[0,1] = "Scotland Championship"
$row213_bg[0,2] = "Scotland Championship"
$row213_bg[0,3] = 45349.69791666666
$row213_bg[0,4] = "Dundee Utd"
$row213_bg[0,5] = "Airdrieonians"
$row213_kv = New-Object 'object[,]' 1,12
$row213_kv[0,0] = 1.4
$row213_kv[0,1] = 4.333
$row213_kv[0,2] = 7
$row213_kv[0,3] = 1.4
$row213_kv[0,4] = 4.333
$row213_kv[0,5] = 7
$row213_kv[0,6] = -1.25
$row213_kv[0,7] = 1.95
$row213_kv[0,8] = 1.9
$row213_kv[0,9] = 2.75
$row213_kv[0,10] = 2
$row213_kv[0,11] = 1.85
$row214_bg = New-Object 'object[,]' 1,6
$row214_bg[0,0] = 6845347
$row214_bg[0,1] = "Scotland Championship"
$row214_bg[0,2] = "Scotland Championship"
$row214_bg[0,3] = 45349.69791666666
$row214_bg[0,4] = "Inverness CT"
$row214_bg[0,5] = "Dunfermline"
$row214_kv = New-Object 'object[,]' 1,12
$row214_kv[0,0] = 2
$row214_kv[0,1] = 3.4
$row214_kv[0,2] = 3.4
$row214_kv[0,3] = 2
$row214_kv[0,4] = 3.4
$row214_kv[0,5] = 3.4
$row214_kv[0,6] = -0.5
$row214_kv[0,7] = 2.05
$row214_kv[0,8] = 1.8
$row214_kv[0,9] = 2.5
$row214_kv[0,10] = 1.975
$row214_kv[0,11] = 1.875
$row215_bg = New-Object 'object[,]' 1,6
$row215_bg[0,0] = 6845348
$row215_bg[0,1] = "Scotland Championship"
$row215_bg[0,2] = "Scotland Championship"
$row215_bg[0,3] = 45349.69791666666
$row215_bg[0,4] = "Queens Park"
$row215_bg[0,5] = "Arbroath"
$row215_kv = New-Object 'object[,]' 1,12
$row215_kv[0,0] = 1.727
$row215_kv[0,1] = 3.75
$row215_kv[0,2] = 4
$row215_kv[0,3] = 1.8
$row215_kv[0,4] = 3.75
$row215_kv[0,5] = 3.8
$row215_kv[0,6] = -0.5
$row215_kv[0,7] = 1.8
$row215_kv[0,8] = 2.05
$row215_kv[0,9] = 2.5
$row215_kv[0,10] = 1.85
$row215_kv[0,11] = 2

$ws.Range("B211:G211").Value = $row212_bg
$ws.Range("K211:V211").Value = $row212_kv
$ws.Range("B212:G212").Value = $row213_bg
$ws.Range("K212:V212").Value = $row213_kv
$ws.Range("B213:G213").Value = $row214_bg
$ws.Range("K213:V213").Value = $row214_kv
$ws.Range("B214:G214").Value = $row215_bg
$ws.Range("K214:V214").Value = $row215_kv
$ws.Range("B215:G215").Value = $row211_bg
$ws.Range("K215:V215").Value = $row211_kv
